# Reverse the order of the comma-separated names/emails in column G
# ("Recorded By") for every data row in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's UsedRange.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    $value = $cell.Value2

    if ($null -ne $value -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $reversedParts = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversedParts += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
